$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.541.02"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "'1.960.51"
$ws.Range("E3").Value = "  +1.36%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'244.58"
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = "  +1.29%  "

$ws.Range("D7").Value = "'58.69"
$ws.Range("E7").Value = "  +2.98%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.374"
$ws.Range("E9").Value = "  +4.28%  "

$ws.Range("D10").Value = "'0.0791"
$ws.Range("E10").Value = "  -5.41%  "

$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("E12").Value = "  +6.00%  "

$ws.Range("D13").Value = "'0.839"
$ws.Range("E13").Value = "  +4.49%  "

$ws.Range("D14").Value = "'2.247.72"
$ws.Range("E14").Value = "  +1.48%  "

$ws.Range("D15").Value = "'21.23"
$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("D16").Value = "'5.29"
$ws.Range("E16").Value = "  +2.90%  "

$ws.Range("D17").Value = "'1.959.89"
$ws.Range("E17").Value = "  +0.94%  "

$ws.Range("D18").Value = "'36.540.25"
$ws.Range("E18").Value = "  +0.73%  "

$ws.Range("D19").Value = "'69.79"
$ws.Range("E19").Value = "  +1.25%  "

$ws.Range("D20").Value = "'0.0₃0849"
$ws.Range("E20").Value = "  -1.38%  "

$ws.Range("D21").Value = "'229.85"
$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("E22").Value = "  +2.09%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'2.46"
$ws.Range("E24").Value = "  +4.95%  "

$ws.Range("E25").Value = "  +3.37%  "

$ws.Range("E26").Value = "  -1.21%  "

$ws.Range("D27").Value = "'0.139"
$ws.Range("E27").Value = "  +6.24%  "

$ws.Range("D28").Value = "'160.76"
$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("D29").Value = "'19.43"
$ws.Range("E29").Value = "  +1.37%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.120"
$ws.Range("E30").Value = "  +2.29%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.21"
$ws.Range("E31").Value = "  +8.38%  "

$ws.Range("D32").Value = "'4.72"
$ws.Range("E32").Value = "  +3.77%  "

$ws.Range("E33").Value = "  -2.31%  "

$ws.Range("E34").Value = "  +5.88%  "

$ws.Range("D35").Value = "'3.51"
$ws.Range("E35").Value = "  +17.60%  "

$ws.Range("D36").Value = "'2.29"
$ws.Range("E36").Value = "  +7.95%  "

$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("E38").Value = "  -0.94%  "

$ws.Range("D39").Value = "'5.47"
$ws.Range("E39").Value = "  -10.03%  "

$ws.Range("D40").Value = "'0.0982"
$ws.Range("E40").Value = "  +1.43%  "

$ws.Range("E41").Value = "  +1.45%  "

$ws.Range("E42").Value = "  +1.66%  "

$ws.Range("D43").Value = "'0.0211"
$ws.Range("E43").Value = "  +1.33%  "

$ws.Range("D44").Value = "'1.372.51"
$ws.Range("E44").Value = "  +3.00%  "

$ws.Range("D45").Value = "'15.78"
$ws.Range("E45").Value = "  +1.59%  "

$ws.Range("D46").Value = "'88.13"
$ws.Range("E46").Value = "  +2.09%  "

$ws.Range("D47").Value = "'1.02"
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("D48").Value = "'7.13"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("D49").Value = "'2.84"
$ws.Range("E49").Value = "  +0.99%  "

$ws.Range("D50").Value = "'2.137.83"
$ws.Range("E50").Value = "  +1.40%  "

$ws.Range("D51").Value = "'44.04"
$ws.Range("E51").Value = "  +0.30%  "
